# Insert a new price record row for "Femacal de La Calera - Zanahoria" above
# the existing row 184, shifting the remaining rows (and the final row that
# falls off the end of the original range) down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 184; everything below (including the former
# last data row, 295) shifts down by one row, so the sheet grows to 296 rows.
$ws.Rows.Item(184).Insert()

# Populate the newly inserted row 184 with the new observation. Most of the
# columns mirror what used to be in row 184 (now shifted to row 185); only
# the date (D), min/max/avg price (K/L/M) and $/Kg price (P) are new values.
$ws.Cells.Item(184, 1).Value  = 3
$ws.Cells.Item(184, 2).Value  = 'Femacal de La Calera'
$ws.Cells.Item(184, 3).Value  = 'Coquimbo'
$ws.Cells.Item(184, 4).Value  = 44606
$ws.Cells.Item(184, 5).Value  = 5
$ws.Cells.Item(184, 6).Value  = 100114013
$ws.Cells.Item(184, 7).Value  = 'Zanahoria'
$ws.Cells.Item(184, 8).Value  = 'Sin especificar'
$ws.Cells.Item(184, 9).Value  = 'Primera'
$ws.Cells.Item(184, 10).Value = 340
$ws.Cells.Item(184, 11).Value = 9500
$ws.Cells.Item(184, 12).Value = 10000
$ws.Cells.Item(184, 13).Value = 9765
$ws.Cells.Item(184, 14).Value = '$/saco 20 kilos'
$ws.Cells.Item(184, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(184, 16).Value = 488
$ws.Cells.Item(184, 17).Value = 20
$ws.Cells.Item(184, 18).Value = 'Hortaliza'
